$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.285.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.197.51'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.33%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.56'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.92%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0769'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.99'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.96'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -11.13%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.530.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.22'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.79%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.195.48'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.30%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.204.77'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.56%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.68'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.26'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '224.73'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.80%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.79'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.52'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.41%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.06'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.60'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -9.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.78'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.42%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.33'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.73%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0692'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.03%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.25'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0951'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.80%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.79%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.902.19'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -9.00%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.95'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.95'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -9.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '71.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.402.51'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '87.05'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.82%  '
